# Auto-generated edit script applying the Ultima_Profits.xlsx value updates
# (current/Leve sale price refresh across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1897.5
$ws.Range("I32").Value = 1726.6666
$ws.Range("K32").Value = 1726.6666
$ws.Range("M32").Value = -1400.6666

$ws.Range("H33").Value = 34483080
$ws.Range("I33").Value = 35714616
$ws.Range("J33").Value = 28
$ws.Range("K33").Value = 35714616
$ws.Range("L33").Value = 28
$ws.Range("M33").Value = -35714387
$ws.Range("N33").Value = -486

$ws.Range("H64").Value = 4276104.5
$ws.Range("I64").Value = 7938911
$ws.Range("J64").Value = 2829.9167
$ws.Range("K64").Value = 7938911
$ws.Range("L64").Value = 2829.9167
$ws.Range("M64").Value = -7938663
$ws.Range("N64").Value = -3325.9167

$ws.Range("H67").Value = 4276104.5
$ws.Range("I67").Value = 7938911
$ws.Range("J67").Value = 2829.9167
$ws.Range("K67").Value = 7938911
$ws.Range("L67").Value = 2829.9167
$ws.Range("M67").Value = -7938053
$ws.Range("N67").Value = -4545.9167

$ws.Range("H74").Value = 3236.0833
$ws.Range("I74").Value = 3166.6667
$ws.Range("J74").Value = 3305.5
$ws.Range("K74").Value = 3166.6667
$ws.Range("L74").Value = 3305.5
$ws.Range("M74").Value = -2230.6667
$ws.Range("N74").Value = -5177.5

$ws.Range("H76").Value = 6151.6
$ws.Range("I76").Value = 4575.8823
$ws.Range("J76").Value = 9500
$ws.Range("K76").Value = 4575.8823
$ws.Range("L76").Value = 9500
$ws.Range("M76").Value = -4260.8823
$ws.Range("N76").Value = -10130

$ws.Range("H77").Value = 3236.0833
$ws.Range("I77").Value = 3166.6667
$ws.Range("J77").Value = 3305.5
$ws.Range("K77").Value = 15833.3335
$ws.Range("L77").Value = 16527.5
$ws.Range("M77").Value = -11153.3335
$ws.Range("N77").Value = -25887.5

$ws.Range("H79").Value = 6151.6
$ws.Range("I79").Value = 4575.8823
$ws.Range("J79").Value = 9500
$ws.Range("K79").Value = 4575.8823
$ws.Range("L79").Value = 9500
$ws.Range("M79").Value = -3483.8823
$ws.Range("N79").Value = -11684

$ws.Range("H88").Value = 5890706
$ws.Range("I88").Value = 10142.857
$ws.Range("J88").Value = 10007100
$ws.Range("K88").Value = 10142.857
$ws.Range("L88").Value = 10007100
$ws.Range("M88").Value = -9736.857
$ws.Range("N88").Value = -10007912

$ws.Range("H91").Value = 5890706
$ws.Range("I91").Value = 10142.857
$ws.Range("J91").Value = 10007100
$ws.Range("K91").Value = 10142.857
$ws.Range("L91").Value = 10007100
$ws.Range("M91").Value = -8738.857
$ws.Range("N91").Value = -10009908

$ws.Range("H132").Value = 5733.8647
$ws.Range("I132").Value = 3244.8928
$ws.Range("J132").Value = 13477.333
$ws.Range("K132").Value = 9734.678400000001
$ws.Range("L132").Value = 40431.999
$ws.Range("M132").Value = -7204.678400000001
$ws.Range("N132").Value = -45491.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2758.6553
$ws.Range("I61").Value = 2717.0454
$ws.Range("J61").Value = 2889.4285
$ws.Range("K61").Value = 2717.0454
$ws.Range("L61").Value = 2889.4285
$ws.Range("M61").Value = -2505.0454
$ws.Range("N61").Value = -3313.4285

$ws.Range("H63").Value = 62502310
$ws.Range("I63").Value = 71430900
$ws.Range("J63").Value = 2200
$ws.Range("K63").Value = 71430900
$ws.Range("L63").Value = 2200
$ws.Range("M63").Value = -71430214
$ws.Range("N63").Value = -3572

$ws.Range("H66").Value = 62502310
$ws.Range("I66").Value = 71430900
$ws.Range("J66").Value = 2200
$ws.Range("K66").Value = 357154500
$ws.Range("L66").Value = 11000
$ws.Range("M66").Value = -357151068
$ws.Range("N66").Value = -17864

$ws.Range("H88").Value = 2571.4285
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 2833.3333
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 2833.3333
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -3645.3333

$ws.Range("H91").Value = 2571.4285
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 2833.3333
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 2833.3333
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -5641.3333

$ws.Range("H136").Value = 2758.6553
$ws.Range("I136").Value = 2717.0454
$ws.Range("J136").Value = 2889.4285
$ws.Range("K136").Value = 8151.1362
$ws.Range("L136").Value = 8668.2855
$ws.Range("M136").Value = -5601.1362
$ws.Range("N136").Value = -13768.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 50181
$ws.Range("J62").Value = 50181
$ws.Range("L62").Value = 50181
$ws.Range("N62").Value = -51553

$ws.Range("H65").Value = 50181
$ws.Range("J65").Value = 50181
$ws.Range("L65").Value = 150543
$ws.Range("N65").Value = -157407

$ws.Range("H105").Value = 4006.6667
$ws.Range("I105").Value = 1440
$ws.Range("J105").Value = 4520
$ws.Range("K105").Value = 1440
$ws.Range("L105").Value = 4520
$ws.Range("M105").Value = 307
$ws.Range("N105").Value = -8014

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 446498.53
$ws.Range("I31").Value = 4649.1
$ws.Range("J31").Value = 1226232.9
$ws.Range("K31").Value = 4649.1
$ws.Range("L31").Value = 1226232.9
$ws.Range("M31").Value = -4354.1
$ws.Range("N31").Value = -1226822.9

$ws.Range("H34").Value = 446498.53
$ws.Range("I34").Value = 4649.1
$ws.Range("J34").Value = 1226232.9
$ws.Range("K34").Value = 4649.1
$ws.Range("L34").Value = 1226232.9
$ws.Range("M34").Value = -4447.1
$ws.Range("N34").Value = -1226636.9

$ws.Range("H111").Value = 44801
$ws.Range("J111").Value = 44801
$ws.Range("L111").Value = 44801
$ws.Range("N111").Value = -52981

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H132").Value = 3839
$ws.Range("I132").Value = 2946.5454
$ws.Range("K132").Value = 8839.636200000001
$ws.Range("M132").Value = -6309.636200000001

$ws.Range("H134").Value = 954340.6
$ws.Range("I134").Value = 1966.238
$ws.Range("J134").Value = 5954306
$ws.Range("K134").Value = 5898.714
$ws.Range("L134").Value = 17862918
$ws.Range("M134").Value = -3363.714
$ws.Range("N134").Value = -17867988

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1364.7013
$ws.Range("I131").Value = 1933.4642
$ws.Range("J131").Value = 1039.6938
$ws.Range("K131").Value = 5800.392599999999
$ws.Range("L131").Value = 3119.0814
$ws.Range("M131").Value = -760.3925999999992
$ws.Range("N131").Value = -13199.0814

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 37204.723
$ws.Range("I70").Value = 49691.23
$ws.Range("J70").Value = 4739.8
$ws.Range("K70").Value = 49691.23
$ws.Range("L70").Value = 4739.8
$ws.Range("M70").Value = -49421.23
$ws.Range("N70").Value = -5279.8

$ws.Range("H73").Value = 37204.723
$ws.Range("I73").Value = 49691.23
$ws.Range("J73").Value = 4739.8
$ws.Range("K73").Value = 49691.23
$ws.Range("L73").Value = 4739.8
$ws.Range("M73").Value = -48755.23
$ws.Range("N73").Value = -6611.8

$ws.Range("H132").Value = 3543.392
$ws.Range("I132").Value = 2750.1924
$ws.Range("J132").Value = 4368.32
$ws.Range("K132").Value = 8250.5772
$ws.Range("L132").Value = 13104.96
$ws.Range("M132").Value = -5720.5772
$ws.Range("N132").Value = -18164.96

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11370333
$ws.Range("I132").Value = 3185.2666
$ws.Range("J132").Value = 35728508
$ws.Range("K132").Value = 9555.799800000001
$ws.Range("L132").Value = 107185524
$ws.Range("M132").Value = -7025.799800000001
$ws.Range("N132").Value = -107190584

$ws.Range("H136").Value = 27032918
$ws.Range("I136").Value = 55557864
$ws.Range("J136").Value = 9285
$ws.Range("K136").Value = 166673592
$ws.Range("L136").Value = 27855
$ws.Range("M136").Value = -166671042
$ws.Range("N136").Value = -32955

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 10833342
$ws.Range("J24").Value = 10833342
$ws.Range("L24").Value = 10833342
$ws.Range("N24").Value = -10833802

$ws.Range("H81").Value = 1730.826
$ws.Range("J81").Value = 2149.0833
$ws.Range("L81").Value = 4298.1666
$ws.Range("N81").Value = -6420.1666

$ws.Range("H84").Value = 1730.826
$ws.Range("J84").Value = 2149.0833
$ws.Range("L84").Value = 21490.833
$ws.Range("N84").Value = -32098.833

$ws.Range("H131").Value = 99800
$ws.Range("J131").Value = 99800
$ws.Range("L131").Value = 99800
$ws.Range("N131").Value = -109880

$ws.Range("H132").Value = 3993.9487
$ws.Range("I132").Value = 4130.971
$ws.Range("J132").Value = 2795
$ws.Range("K132").Value = 12392.913
$ws.Range("L132").Value = 8385
$ws.Range("M132").Value = -9862.912999999999
$ws.Range("N132").Value = -13445

$ws.Range("H136").Value = 4348968
$ws.Range("I136").Value = 5001088
$ws.Range("J136").Value = 1501.6666
$ws.Range("K136").Value = 15003264
$ws.Range("L136").Value = 4504.9998
$ws.Range("M136").Value = -15000714
$ws.Range("N136").Value = -9604.9998
